$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: populate with the same pattern used by rows 2 & 3 ---
$ws.Range("A4").Value = "run"
$ws.Range("B4").Value = "DPLKKLM052-001"
$ws.Range("C4").Value = "Klaim - Transaksi - Validasi Klaim Pasca Kerja "
$ws.Range("D4").Value = "Klaim Pasca Kerja "
$ws.Range("E4").Value = "Validasi Klaim bisa dilakukan dengan baik. Dalam perhitungan nominal klaim, dikenakan biaya administrasi dan fee yang disesuaikan dengan ketentuan PKS"
$ws.Range("G4").Value = 48968
$ws.Range("H4").Value = "bni1234"
$ws.Range("I4").Value = "Klaim"
$ws.Range("J4").Formula = "'Transaksi"
$ws.Range("K4").Value = "Validasi Klaim Pasca Kerja "
$ws.Range("N4").Formula = "'0000000045"
$ws.Range("F4").Formula = @'
= "Username : "&G4&",
Password : bni1234,
ID Peserta :  "&N4
'@
$ws.Rows.Item(4).RowHeight = 75

# --- Row 3: drop the stray A3 cell (its content moved to A4) ---
$ws.Range("A3").Clear()

# --- Selection / view: move to A3, no horizontal scroll needed anymore ---
$ws.Range("A3").Select()
